$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at column AZ ("Mean"), the average of the 50
# "Run" columns (B:AY). A new run (Run 50) was added as its own column,
# so insert a fresh column at AZ - this pushes the existing "Mean" column
# (with all its data and formatting) one slot to the right, into BA.
$ws.Columns("AZ").Insert()

# Fill in the new "Run 50" header in AZ1, matching the bold/bordered
# header style used by the rest of row 1.
$ws.Range("AZ1").Value() = "Run 50"
$ws.Range("AZ1").Font.Bold = $ws.Range("AY1").Font.Bold
$ws.Range("AZ1").Borders.LineStyle = 1
$ws.Range("AZ1").HorizontalAlignment = $ws.Range("AY1").HorizontalAlignment
$ws.Range("AZ1").VerticalAlignment = $ws.Range("AY1").VerticalAlignment

# New per-row "Run 50" result, and the recomputed "Mean" (now averaging
# B:AZ, i.e. the original 50 runs plus the new Run 50 column), for each
# of the 13 data rows.
$run50 = 8.048758810000001
$newMean = 6.54283568

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value() = $run50
    $ws.Cells.Item($r, 53).Value() = $newMean
}

Write-Output "Added Run 50 column and refreshed Mean column"
